$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from the default "Sheet1" to "Query Results".
$ws.Name = "Query Results"

# The source workbook was saved with a stale selection left over from
# browsing the data (A3:XFD91500). Reset the selection back to the
# top-left cell, which is the state a freshly generated export would be
# in.
$ws.Range("A1").Select() | Out-Null
